# Update gh-pages to output generated at 456a3b4
# Sheet "展览" (index 1) and sheet "全部类型" (index 4) both get refreshed
# "想去人数" (F) and "最低票价" (G) figures for the same events.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 16722
$ws1.Range("G2").Value = "暂时售罄"
$ws1.Range("F3").Value = 356
$ws1.Range("F6").Value = 716
$ws1.Range("F7").Value = 1769
$ws1.Range("F8").Value = 168

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 16723
$ws4.Range("G2").Value = "暂时售罄"
$ws4.Range("F3").Value = 356
$ws4.Range("F8").Value = 716
$ws4.Range("F9").Value = 1769
$ws4.Range("F11").Value = 168
